# Update "想去人数" (interest count) values on the "展览" and "全部类型" sheets
# to reflect the refreshed data pull (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 8536
    $ws.Range("F4").Value = 385
    $ws.Range("F5").Value = 26
}
